# The "reviews_count" column (column E) was removed from the sheet.
# Deleting the entire column shifts every column to its right
# (reviews_average, latitude, longitude, is_permanently_closed,
# gmaps_link, latest_review_date) one position to the left,
# which matches the target diff (columns F:K -> E:J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(5).Delete()
